$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-04-04 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-05 Wednesday", 2) | Out-Null

# Update the 20x5 multiplication table, cell by cell, in row-major order
$t = $d.Tables.Item(1)
$rows = 20
$cols = 5

$t.Cell(1, 1).Range.Text = "96×31="
$t.Cell(1, 2).Range.Text = "76×39="
$t.Cell(1, 3).Range.Text = "88×22="
$t.Cell(1, 4).Range.Text = "66×30="
$t.Cell(1, 5).Range.Text = "71×18="
$t.Cell(2, 1).Range.Text = "81×53="
$t.Cell(2, 2).Range.Text = "47×46="
$t.Cell(2, 3).Range.Text = "42×48="
$t.Cell(2, 4).Range.Text = "95×35="
$t.Cell(2, 5).Range.Text = "21×88="
$t.Cell(3, 1).Range.Text = "96×23="
$t.Cell(3, 2).Range.Text = "40×63="
$t.Cell(3, 3).Range.Text = "86×13="
$t.Cell(3, 4).Range.Text = "51×67="
$t.Cell(3, 5).Range.Text = "66×82="
$t.Cell(4, 1).Range.Text = "48×49="
$t.Cell(4, 2).Range.Text = "28×59="
$t.Cell(4, 3).Range.Text = "55×43="
$t.Cell(4, 4).Range.Text = "41×32="
$t.Cell(4, 5).Range.Text = "42×20="
$t.Cell(5, 1).Range.Text = "63×52="
$t.Cell(5, 2).Range.Text = "30×44="
$t.Cell(5, 3).Range.Text = "16×58="
$t.Cell(5, 4).Range.Text = "18×56="
$t.Cell(5, 5).Range.Text = "23×61="
$t.Cell(6, 1).Range.Text = "67×45="
$t.Cell(6, 2).Range.Text = "85×10="
$t.Cell(6, 3).Range.Text = "51×67="
$t.Cell(6, 4).Range.Text = "96×28="
$t.Cell(6, 5).Range.Text = "45×56="
$t.Cell(7, 1).Range.Text = "13×26="
$t.Cell(7, 2).Range.Text = "32×24="
$t.Cell(7, 3).Range.Text = "46×48="
$t.Cell(7, 4).Range.Text = "92×97="
$t.Cell(7, 5).Range.Text = "58×15="
$t.Cell(8, 1).Range.Text = "96×27="
$t.Cell(8, 2).Range.Text = "73×38="
$t.Cell(8, 3).Range.Text = "68×79="
$t.Cell(8, 4).Range.Text = "51×64="
$t.Cell(8, 5).Range.Text = "50×60="
$t.Cell(9, 1).Range.Text = "87×79="
$t.Cell(9, 2).Range.Text = "46×30="
$t.Cell(9, 3).Range.Text = "37×47="
$t.Cell(9, 4).Range.Text = "79×12="
$t.Cell(9, 5).Range.Text = "69×32="
$t.Cell(10, 1).Range.Text = "41×55="
$t.Cell(10, 2).Range.Text = "19×16="
$t.Cell(10, 3).Range.Text = "54×34="
$t.Cell(10, 4).Range.Text = "53×36="
$t.Cell(10, 5).Range.Text = "96×64="
$t.Cell(11, 1).Range.Text = "47×17="
$t.Cell(11, 2).Range.Text = "44×68="
$t.Cell(11, 3).Range.Text = "81×27="
$t.Cell(11, 4).Range.Text = "81×30="
$t.Cell(11, 5).Range.Text = "79×22="
$t.Cell(12, 1).Range.Text = "13×12="
$t.Cell(12, 2).Range.Text = "24×21="
$t.Cell(12, 3).Range.Text = "90×74="
$t.Cell(12, 4).Range.Text = "59×73="
$t.Cell(12, 5).Range.Text = "26×61="
$t.Cell(13, 1).Range.Text = "90×20="
$t.Cell(13, 2).Range.Text = "46×52="
$t.Cell(13, 3).Range.Text = "54×18="
$t.Cell(13, 4).Range.Text = "90×25="
$t.Cell(13, 5).Range.Text = "95×77="
$t.Cell(14, 1).Range.Text = "65×96="
$t.Cell(14, 2).Range.Text = "56×97="
$t.Cell(14, 3).Range.Text = "41×54="
$t.Cell(14, 4).Range.Text = "56×81="
$t.Cell(14, 5).Range.Text = "31×58="
$t.Cell(15, 1).Range.Text = "80×25="
$t.Cell(15, 2).Range.Text = "59×23="
$t.Cell(15, 3).Range.Text = "22×83="
$t.Cell(15, 4).Range.Text = "30×31="
$t.Cell(15, 5).Range.Text = "88×65="
$t.Cell(16, 1).Range.Text = "11×69="
$t.Cell(16, 2).Range.Text = "92×97="
$t.Cell(16, 3).Range.Text = "90×61="
$t.Cell(16, 4).Range.Text = "85×85="
$t.Cell(16, 5).Range.Text = "22×29="
$t.Cell(17, 1).Range.Text = "49×13="
$t.Cell(17, 2).Range.Text = "64×11="
$t.Cell(17, 3).Range.Text = "75×54="
$t.Cell(17, 4).Range.Text = "49×12="
$t.Cell(17, 5).Range.Text = "50×44="
$t.Cell(18, 1).Range.Text = "39×75="
$t.Cell(18, 2).Range.Text = "78×29="
$t.Cell(18, 3).Range.Text = "45×95="
$t.Cell(18, 4).Range.Text = "13×73="
$t.Cell(18, 5).Range.Text = "58×38="
$t.Cell(19, 1).Range.Text = "17×50="
$t.Cell(19, 2).Range.Text = "88×81="
$t.Cell(19, 3).Range.Text = "59×96="
$t.Cell(19, 4).Range.Text = "20×76="
$t.Cell(19, 5).Range.Text = "58×39="
$t.Cell(20, 1).Range.Text = "96×92="
$t.Cell(20, 2).Range.Text = "56×23="
$t.Cell(20, 3).Range.Text = "98×68="
$t.Cell(20, 4).Range.Text = "34×76="
$t.Cell(20, 5).Range.Text = "76×49="
